$d = $word.ActiveDocument

# Replace "present" with "February 2020" in the Softwarium line
$d.Content.Find.Execute("Web-developer Softwarium in Kyiv April 2019 - present", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Web-developer Softwarium in Kyiv April 2019 - February 2020", 2)
